$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataframe's rows/columns are being swapped (transposed):
#   - Old layout: column A (rows 2-4) held the field labels, column B (rows 1-4) held the values
#   - New layout: row 1 (cols B-D) holds the field labels, row 2 (cols A-D) holds the values

# Step 1: stash A2 ("VideoName" label, with its header style) in a scratch cell
# so it isn't lost when B1 is copied into A2.
$ws.Range("A2").Copy($ws.Range("Z1"))

# Step 2: move the numeric id that was in B1 down into A2 (keeps header style).
$ws.Range("B1").Copy($ws.Range("A2"))

# Step 3: move the stashed "VideoName" label (with style) into its new home, B1.
$ws.Range("Z1").Copy($ws.Range("B1"))

# Step 4: move the "Views" label (with style) from A3 into C1.
$ws.Range("A3").Copy($ws.Range("C1"))

# Step 5: move the "Date" label (with style) from A4 into D1.
$ws.Range("A4").Copy($ws.Range("D1"))

# Step 6: update the video name value in place (B2 keeps its position).
$ws.Range("B2").Value = "Brushed DC Motors and How to Drive Them"

# Step 7: the view-count value moves from B3 to C2, with a new number.
$ws.Range("C2").Value = 3227096

# Step 8: the date value moves from B4 to D2, with a new date string.
$ws.Range("D2").Value = "25 Jan 2019"

# Step 9: clear out the now-unused cells from the old layout.
$ws.Range("A1").ClearContents()
$ws.Range("A3").Clear()
$ws.Range("A4").Clear()
$ws.Range("B3").Clear()
$ws.Range("B4").Clear()
$ws.Range("Z1").Clear()
